# Updated symbol list on Sun Jan 22 10:15:48 UTC 2023 with GitHub Actions
# Refresh crypto price/volume/hour snapshot data on the "cryptos" worksheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry below is a cell address together with its refreshed text value.
# All target cells (Price, Volume(1h) and Hora columns) hold plain text in the
# original workbook, so we force a text number format before writing the new
# value to avoid Excel silently re-interpreting numeric-looking strings
# (e.g. "301.82" or "10") as numbers.
$updates = @(
    @{ Cell = "D2"; Value = '301.82' },
    @{ Cell = "E2"; Value = '-1.81%' },
    @{ Cell = "G2"; Value = '10' },
    @{ Cell = "D3"; Value = '37.38' },
    @{ Cell = "E3"; Value = '6.74%' },
    @{ Cell = "G3"; Value = '10' },
    @{ Cell = "E4"; Value = '-3.84%' },
    @{ Cell = "G4"; Value = '10' },
    @{ Cell = "D5"; Value = '0.07804' },
    @{ Cell = "E5"; Value = '-1.20%' },
    @{ Cell = "G5"; Value = '10' },
    @{ Cell = "D6"; Value = '2.197' },
    @{ Cell = "E6"; Value = '-6.78%' },
    @{ Cell = "G6"; Value = '10' },
    @{ Cell = "D7"; Value = '8.032' },
    @{ Cell = "E7"; Value = '0.15%' },
    @{ Cell = "G7"; Value = '10' },
    @{ Cell = "D8"; Value = '4.037' },
    @{ Cell = "E8"; Value = '1.59%' },
    @{ Cell = "G8"; Value = '10' },
    @{ Cell = "D9"; Value = '0.9149' },
    @{ Cell = "E9"; Value = '-1.61%' },
    @{ Cell = "G9"; Value = '10' },
    @{ Cell = "D10"; Value = '0.09663' },
    @{ Cell = "E10"; Value = '-5.20%' },
    @{ Cell = "G10"; Value = '10' },
    @{ Cell = "D11"; Value = '0.1883' },
    @{ Cell = "E11"; Value = '2.24%' },
    @{ Cell = "G11"; Value = '10' },
    @{ Cell = "D12"; Value = '0.08730' },
    @{ Cell = "E12"; Value = '0.57%' },
    @{ Cell = "G12"; Value = '10' },
    @{ Cell = "D13"; Value = '0.03526' },
    @{ Cell = "E13"; Value = '4.88%' },
    @{ Cell = "G13"; Value = '10' },
    @{ Cell = "D14"; Value = '0.09953' },
    @{ Cell = "E14"; Value = '0.20%' },
    @{ Cell = "G14"; Value = '10' },
    @{ Cell = "D15"; Value = '0.001494' },
    @{ Cell = "E15"; Value = '0.84%' },
    @{ Cell = "G15"; Value = '10' },
    @{ Cell = "D16"; Value = '0.005693' },
    @{ Cell = "E16"; Value = '0.35%' },
    @{ Cell = "G16"; Value = '10' },
    @{ Cell = "D17"; Value = '3.463' },
    @{ Cell = "E17"; Value = '-0.66%' },
    @{ Cell = "G17"; Value = '10' },
    @{ Cell = "D18"; Value = '2.378' },
    @{ Cell = "E18"; Value = '10.97%' },
    @{ Cell = "G18"; Value = '10' },
    @{ Cell = "E19"; Value = '1.88%' },
    @{ Cell = "G19"; Value = '10' },
    @{ Cell = "D20"; Value = '0.1277' },
    @{ Cell = "E20"; Value = '-2.18%' },
    @{ Cell = "G20"; Value = '10' },
    @{ Cell = "D21"; Value = '4.773' },
    @{ Cell = "E21"; Value = '4.85%' },
    @{ Cell = "G21"; Value = '10' },
    @{ Cell = "D22"; Value = '0.2296' },
    @{ Cell = "E22"; Value = '0.25%' },
    @{ Cell = "G22"; Value = '10' },
    @{ Cell = "D23"; Value = '0.04629' },
    @{ Cell = "E23"; Value = '1.35%' },
    @{ Cell = "G23"; Value = '10' },
    @{ Cell = "E24"; Value = '1.50%' },
    @{ Cell = "G24"; Value = '10' },
    @{ Cell = "D25"; Value = '0.004788' },
    @{ Cell = "E25"; Value = '7.11%' },
    @{ Cell = "G25"; Value = '10' },
    @{ Cell = "E26"; Value = '-6.79%' },
    @{ Cell = "G26"; Value = '10' },
    @{ Cell = "E27"; Value = '40.10%' },
    @{ Cell = "G27"; Value = '10' },
    @{ Cell = "G28"; Value = '10' },
    @{ Cell = "G29"; Value = '10' },
    @{ Cell = "G30"; Value = '10' },
    @{ Cell = "G31"; Value = '10' },
    @{ Cell = "G32"; Value = '10' },
    @{ Cell = "G33"; Value = '10' },
    @{ Cell = "G34"; Value = '10' },
    @{ Cell = "G35"; Value = '10' },
    @{ Cell = "G36"; Value = '10' },
    @{ Cell = "G37"; Value = '10' },
    @{ Cell = "G38"; Value = '10' },
    @{ Cell = "D39"; Value = '0.01756' },
    @{ Cell = "E39"; Value = '-1.55%' },
    @{ Cell = "G39"; Value = '10' },
    @{ Cell = "D40"; Value = '0.04732' },
    @{ Cell = "E40"; Value = '-1.59%' },
    @{ Cell = "G40"; Value = '10' },
    @{ Cell = "D41"; Value = '0.008079' },
    @{ Cell = "E41"; Value = '3.84%' },
    @{ Cell = "G41"; Value = '10' },
    @{ Cell = "D42"; Value = '0.1390' },
    @{ Cell = "E42"; Value = '-1.88%' },
    @{ Cell = "G42"; Value = '10' },
    @{ Cell = "D43"; Value = '0.007687' },
    @{ Cell = "E43"; Value = '9.02%' },
    @{ Cell = "G43"; Value = '10' },
    @{ Cell = "D44"; Value = '0.002222' },
    @{ Cell = "E44"; Value = '1.00%' },
    @{ Cell = "G44"; Value = '10' },
    @{ Cell = "D45"; Value = '0.01039' },
    @{ Cell = "E45"; Value = '9.10%' },
    @{ Cell = "G45"; Value = '10' },
    @{ Cell = "D46"; Value = '0.00006056' },
    @{ Cell = "E46"; Value = '1.46%' },
    @{ Cell = "G46"; Value = '10' },
    @{ Cell = "E47"; Value = '0.67%' },
    @{ Cell = "G47"; Value = '10' },
    @{ Cell = "D48"; Value = '7.840' },
    @{ Cell = "E48"; Value = '186.94%' },
    @{ Cell = "G48"; Value = '10' },
    @{ Cell = "E49"; Value = '0.32%' },
    @{ Cell = "G49"; Value = '10' },
    @{ Cell = "D50"; Value = '0.00002102' },
    @{ Cell = "E50"; Value = '0.67%' },
    @{ Cell = "G50"; Value = '10' },
    @{ Cell = "D51"; Value = '0.0002002' },
    @{ Cell = "E51"; Value = '0.67%' },
    @{ Cell = "G51"; Value = '10' }
)

foreach ($update in $updates) {
    $cell = $ws.Range($update.Cell)
    $cell.NumberFormat = "@"
    $cell.Value = $update.Value
}
